# TestData.xlsx — "Get Team" test-case row gets parameterised: the response
# map now captures every account field (account_id, locks, paid flags,
# quotas, role_code) for both team members instead of just email/quota.
# Header column widths are re-tuned for the new (shorter per-column) text and
# the key cell (F2) gets centred/wrapped formatting with a taller header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update ResponseMapValues (G5) then ResponseMapKeys (F5) for the
#     "Get Team" row so the two new shared-string entries land in the same
#     slot order as the authored workbook. ---
$ws.Range("G5").Value = "Sofbang Team;7953efb9fd8a30b65b34b03b86bb84de374b4e4f;vivek.ahuja@sofbang.com;false;true;false;null;null;5000;a;538da81b5f9cc2d0faef17f4efb28ee4b95a8c42;puneet.gandhi@sofbang.com;false;true;false;null;null;5000;a;"
$ws.Range("F5").Value = "team.name;team.accounts[0].account_id;team.accounts[0].email_address;team.accounts[0].is_locked;team.accounts[0].is_paid_hs;team.accounts[0].is_paid_hf;team.accounts[0].quotas.templates_left;team.accounts[0].quotas.documents_left;team.accounts[0].quotas.api_signature_requests_left;team.accounts[0].role_code;team.accounts[1].account_id;team.accounts[1].email_address;team.accounts[1].is_locked;team.accounts[1].is_paid_hs;team.accounts[1].is_paid_hf;team.accounts[1].quotas.templates_left;team.accounts[1].quotas.documents_left;team.accounts[1].quotas.api_signature_requests_left;team.accounts[1].role_code;"

# --- New formatting (centered / middle / wrap) for the ResponseMapKeys
#     header-ish cell in row 2 (F2), plus a taller row to fit the wrapped
#     text. ---
$f2 = $ws.Range("F2")
$f2.HorizontalAlignment = -4108   # xlCenter
$f2.VerticalAlignment = -4108     # xlCenter
$f2.WrapText = $true
$ws.Rows.Item(2).RowHeight = 45

# --- F5 picks up a distinct (but visually equivalent / default-aligned)
#     style slot of its own. ---
$f5 = $ws.Range("F5")
$f5.WrapText = $false

# --- Columns F/G are narrower now that the text wraps / reflows. ---
$ws.Columns.Item(6).ColumnWidth = 38.16666666666666
$ws.Columns.Item(7).ColumnWidth = 32

# --- Move the selection/cursor to G7, matching the saved view state. ---
$ws.Range("G7").Select()
